$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.229.18'
$ws.Range("E2").Value = '  -3.34%  '
$ws.Range("D3").Value = '2.840.45'
$ws.Range("E3").Value = '  -3.93%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''505.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.47%  '
$ws.Range("D6").Value = '''135.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.11%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '''0.529'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.25%  '
$ws.Range("D9").Value = '2.838.44'
$ws.Range("E9").Value = '  -4.27%  '
$ws.Range("D10").Value = '''0.104'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.61%  '
$ws.Range("D11").Value = '''5.94'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.70%  '
$ws.Range("D12").Value = '''0.350'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.70%  '
$ws.Range("D13").Value = '3.331.56'
$ws.Range("E13").Value = '  -4.28%  '
$ws.Range("E14").Value = '  +1.56%  '
$ws.Range("D15").Value = '59.327.27'
$ws.Range("E15").Value = '  -3.32%  '
$ws.Range("D16").Value = '''21.90'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.91%  '
$ws.Range("D17").Value = '2.847.57'
$ws.Range("E17").Value = '  -4.17%  '
$ws.Range("D18").Value = '''0.0000136'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.89%  '
$ws.Range("D19").Value = '''4.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.39%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''353.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.68%  '
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = '''11.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.34%  '
$ws.Range("D22").Value = '''6.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.68%  '
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D25").Value = '''63.17'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.20%  '
$ws.Range("D26").Value = '''0.430'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.33%  '
$ws.Range("D27").Value = '''0.172'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.65%  '
$ws.Range("D28").Value = '''1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = '''7.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.00%  '
$ws.Range("D30").Value = '0.0₃0815'
$ws.Range("E30").Value = '  -9.42%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  -4.98%  '
$ws.Range("D33").Value = '''19.08'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.34%  '
$ws.Range("D34").Value = '''150.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.40%  '
$ws.Range("D35").Value = '''4.19'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.74%  '
$ws.Range("D36").Value = '''5.38'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.08%  '
$ws.Range("D37").Value = '''0.908'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -12.41%  '
$ws.Range("E38").Value = '  -8.13%  '
$ws.Range("D39").Value = '''36.57'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.45%  '
$ws.Range("D40").Value = '2.222.13'
$ws.Range("E40").Value = '  -7.15%  '
$ws.Range("D41").Value = '''0.631'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.73%  '
$ws.Range("D42").Value = '''3.55'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.78%  '
$ws.Range("D43").Value = '''1.38'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -8.86%  '
$ws.Range("D44").Value = '''0.0562'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.38%  '
$ws.Range("D45").Value = '''0.999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").Value = '''19.53'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -9.63%  '
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("D48").Value = '''0.0227'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.60%  '
$ws.Range("D49").Value = '''0.0890'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.25%  '
$ws.Range("D50").Value = '''4.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -11.20%  '
$ws.Range("D51").Value = '''17.74'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.91%  '
